$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Whole")

# Year label added above the Jan..Dec header row
$ws.Range("A1").Value = 2021

# Row 3 "Feed Mass": monthly breakdown, all zero except Apr (100)
$ws.Range("B3:M3").Value = 0
$ws.Range("F3").Value = 100

# Row 4 "Feed Price": monthly breakdown, all zero except Apr (644)
$ws.Range("B4:M4").Value = 0
$ws.Range("F4").Value = 644

# Row 5 "Misc": monthly breakdown, all zero
$ws.Range("B5:M5").Value = 0

# touch C6 so it is materialised as a (blank) formatted cell, matching the source
$ws.Range("C6").NumberFormat = "General"

# New summary rows
$ws.Range("A6").Value = "Average Age"
$ws.Range("A7").Value = "Feed per Pig"
$ws.Range("A8").Value = "Feed per Pig per age"

# Column A is widened to fit the new row labels
$ws.Columns.Item(1).ColumnWidth = 19.75

# Leave the selection on F6, matching the edited workbook
$ws.Range("F6").Select() | Out-Null
